$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 186
$ws.Range("F4").Value = 404
$ws.Range("F6").Value = 5200
$ws.Range("F7").Value = 435
$ws.Range("F8").Value = 613
$ws.Range("F9").Value = 901
$ws.Range("F10").Value = 811
$ws.Range("F11").Value = 71
$ws.Range("F13").Value = 556
$ws.Range("F14").Value = 9
$ws.Range("F17").Value = 1726
$ws.Range("F18").Value = 1443
$ws.Range("F19").Value = 806
$ws.Range("F21").Value = 183
$ws.Range("F22").Value = 298
$ws.Range("F23").Value = 502
$ws.Range("F24").Value = 128
$ws.Range("F28").Value = 2453
$ws.Range("F30").Value = 93
$ws.Range("F31").Value = 37
$ws.Range("F32").Value = 84
$ws.Range("F34").Value = 244
$ws.Range("F39").Value = 272
$ws.Range("F40").Value = 628
$ws.Range("F42").Value = 44
$ws.Range("F43").Value = 40

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 148
$ws.Range("F6").Value = 104

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 186
$ws.Range("F7").Value = 5200
$ws.Range("F8").Value = 435
$ws.Range("F9").Value = 613
$ws.Range("F11").Value = 148
$ws.Range("F12").Value = 901
$ws.Range("F13").Value = 811
$ws.Range("F15").Value = 71
$ws.Range("F17").Value = 556
$ws.Range("F18").Value = 9
$ws.Range("F22").Value = 1726
$ws.Range("F23").Value = 1443
$ws.Range("F24").Value = 806
$ws.Range("F26").Value = 183
$ws.Range("F27").Value = 298
$ws.Range("F29").Value = 502
$ws.Range("F30").Value = 128
$ws.Range("F33").Value = 2454
$ws.Range("F35").Value = 93
$ws.Range("F36").Value = 84
$ws.Range("F38").Value = 244
$ws.Range("F42").Value = 272
$ws.Range("F43").Value = 628
$ws.Range("F45").Value = 40

